$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.308.94'
$ws.Range('E2').Value = '  +4.39%  '
$ws.Range('D3').Value = '3.643.07'
$ws.Range('E3').Value = '  +7.22%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '594.47'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').Value = '181.42'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').Value = '3.631.40'
$ws.Range('E7').Value = '  +7.13%  '
$ws.Range('E8').Value = '  +2.00%  '
$ws.Range('E10').Value = '  +4.09%  '
$ws.Range('D11').Value = '0.605'
$ws.Range('E11').Value = '  +1.96%  '
$ws.Range('D12').Value = '50.07'
$ws.Range('E12').Value = '  +2.86%  '
$ws.Range('E13').Value = '  +1.85%  '
$ws.Range('D14').Value = '699.61'
$ws.Range('E14').Value = '  +2.56%  '
$ws.Range('D15').Value = '4.229.92'
$ws.Range('E15').Value = '  +7.43%  '
$ws.Range('D16').Value = '8.96'
$ws.Range('E16').Value = '  +3.42%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.701.59'
$ws.Range('E17').Value = '  +9.12%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '72.431.02'
$ws.Range('E18').Value = '  +4.36%  '
$ws.Range('D20').Value = '18.49'
$ws.Range('E20').Value = '  +3.93%  '
$ws.Range('D21').Value = '11.62'
$ws.Range('E21').Value = '  +2.53%  '
$ws.Range('E22').Value = '  +2.84%  '
$ws.Range('D23').Value = '5.86'
$ws.Range('E23').Value = '  +8.61%  '
$ws.Range('D24').Value = '17.98'
$ws.Range('E24').Value = '  +4.47%  '
$ws.Range('D25').Value = '103.87'
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('E26').Value = '  +2.68%  '
$ws.Range('E27').Value = '  +4.65%  '
$ws.Range('D28').Value = '9.96'
$ws.Range('E28').Value = '  +2.22%  '
$ws.Range('D29').Value = '35.21'
$ws.Range('E29').Value = '  +3.46%  '
$ws.Range('D30').Value = '9.14'
$ws.Range('E30').Value = '  +3.73%  '
$ws.Range('E31').Value = '  +6.28%  '
$ws.Range('D32').Value = '4.18'
$ws.Range('E32').Value = '  +15.80%  '
$ws.Range('D33').Value = '585.77'
$ws.Range('E33').Value = '  +5.12%  '
$ws.Range('D34').Value = '11.29'
$ws.Range('E34').Value = '  +1.17%  '
$ws.Range('D35').Value = '0.109'
$ws.Range('E35').Value = '  +1.87%  '
$ws.Range('D36').Value = '59.82'
$ws.Range('E36').Value = '  +2.17%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').Value = '3.638.86'
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('E39').Value = '  +2.61%  '
$ws.Range('D40').Value = '0.0₃0776'
$ws.Range('E40').Value = '  +7.62%  '
$ws.Range('D41').Value = '35.81'
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  +5.68%  '
$ws.Range('E43').Value = '  +3.40%  '
$ws.Range('D44').Value = '0.0459'
$ws.Range('E44').Value = '  +8.06%  '
$ws.Range('E45').Value = '  +3.17%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').Value = '2.86'
$ws.Range('E46').Value = '  +6.33%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '1.47'
$ws.Range('E47').Value = '  +5.51%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '0.132'
$ws.Range('E48').Value = '  +1.99%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '133.79'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D51').Value = '3.00'
$ws.Range('E51').Value = '  +14.40%  '
